$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.083
$ws.Range("A6").Value = -22.056
$ws.Range("A7").Value = -21.191
$ws.Range("B7").Value = 6.765000000000001
$ws.Range("A8").Value = -21.623
$ws.Range("B11").Value = 6.017
$ws.Range("B12").Value = 5.458
$ws.Range("E12").Value = 17.641
$ws.Range("E13").Value = 16.561
$ws.Range("E14").Value = 17.102
$ws.Range("B15").Value = 5.081
$ws.Range("A16").Value = -21.397
$ws.Range("E16").Value = 16.784
$ws.Range("E19").Value = 16.694
$ws.Range("A20").Value = -21.795
$ws.Range("B20").Value = 5.353999999999999
$ws.Range("E20").Value = 16.347
$ws.Range("A21").Value = -21.26
$ws.Range("B21").Value = 7.531999999999999
$ws.Range("B22").Value = 6.980000000000001
$ws.Range("E22").Value = 16.627
$ws.Range("B23").Value = 7.007
$ws.Range("A28").Value = -21.778
$ws.Range("A29").Value = -21.764
$ws.Range("B29").Value = 5.702
$ws.Range("A30").Value = -21.79499999999999
$ws.Range("A32").Value = -21.705
$ws.Range("B34").Value = 8.059000000000001
$ws.Range("E36").Value = 16.748
$ws.Range("A40").Value = -20.605
$ws.Range("B42").Value = 7.047000000000001
$ws.Range("B43").Value = 5.615
$ws.Range("E43").Value = 17.295
$ws.Range("B44").Value = 5.374000000000001
$ws.Range("B45").Value = 5.25
$ws.Range("A46").Value = -21.017
$ws.Range("B46").Value = 7.31
$ws.Range("E46").Value = 16.695
$ws.Range("B50").Value = 5.393
$ws.Range("E50").Value = 16.497
$ws.Range("A51").Value = -21.26
$ws.Range("B51").Value = 7.411000000000001
$ws.Range("A52").Value = -21.493
$ws.Range("A57").Value = -21.849
$ws.Range("B57").Value = 6.078
$ws.Range("A59").Value = -22.425
$ws.Range("A62").Value = -21.854
$ws.Range("B65").Value = 5.306
$ws.Range("A66").Value = -21.47
$ws.Range("B66").Value = 5.709
$ws.Range("B67").Value = 6.014
$ws.Range("A73").Value = -20.783
$ws.Range("A74").Value = -20.947
$ws.Range("E76").Value = 16.659
$ws.Range("A77").Value = -21.505
$ws.Range("B79").Value = 5.680999999999999
$ws.Range("B84").Value = 5.781000000000001
$ws.Range("B87").Value = 4.447000000000001
$ws.Range("A92").Value = -21.477
$ws.Range("B92").Value = 5.681
$ws.Range("E95").Value = 17.072
$ws.Range("B97").Value = 5.677
$ws.Range("E97").Value = 16.922
$ws.Range("E99").Value = 16.732
$ws.Range("A100").Value = -21.481
